# Insert the new FTT-IH-CHI worksheet after FTT-H, matching the
# "Fixed a few bugs which lets CHI run" commit: adds a brand new sheet of
# variable metadata for the FTT-IH-CHI (industrial heat, chemicals & other
# industry) sub-model between "FTT-H" and "FTT-S".

$wb = $excel.ActiveWorkbook

$wsH = $wb.Worksheets.Item("FTT-H")

# New sheet goes right after FTT-H, before FTT-S/FTT-Fr/Time_Horizons.
$newSheet = $wb.Worksheets.Add($null, $wsH)
$newSheet.Name = "FTT-IH-CHI"

$arr = New-Object 'object[,]' 23,9
$arr[0,0] = "Variable name"
$arr[0,1] = "Read in?"
$arr[0,2] = "Code"
$arr[0,3] = "Description"
$arr[0,4] = "RowDim"
$arr[0,5] = "ColDim"
$arr[0,6] = "3DDim"
$arr[0,7] = "Conversion?"
$arr[0,8] = "Scenario"
$arr[1,0] = "IUD1"
$arr[1,1] = 1
$arr[1,2] = 0
$arr[1,3] = "FTT-IH-CHI useful energy demand (GWh)"
$arr[1,4] = "ITTI"
$arr[1,5] = "TIME"
$arr[1,6] = "RSHORTTI"
$arr[1,7] = 0
$arr[1,8] = "S0"
$arr[2,0] = "ISC1"
$arr[2,1] = 1
$arr[2,2] = 0
$arr[2,3] = "FTT-IH-CHI market share caps"
$arr[2,4] = "ITTI"
$arr[2,5] = "TIME"
$arr[2,6] = "RSHORTTI"
$arr[2,7] = 0
$arr[2,8] = "S0"
$arr[3,0] = "IWW1"
$arr[3,1] = 1
$arr[3,2] = 0
$arr[3,3] = "FTT-IH-CHI cumulative capacities (EU28)"
$arr[3,4] = "ITTI"
$arr[3,5] = 0
$arr[3,6] = 0
$arr[3,7] = "TIME"
$arr[3,8] = "S0"
$arr[4,0] = "IWI1"
$arr[4,1] = 1
$arr[4,2] = 0
$arr[4,3] = "FTT-IH-CHI yearly capacity additions"
$arr[4,4] = "ITTI"
$arr[4,5] = "TIME"
$arr[4,6] = "RSHORTTI"
$arr[4,7] = 0
$arr[4,8] = "S0"
$arr[5,0] = "IWK1"
$arr[5,1] = 1
$arr[5,2] = 0
$arr[5,3] = "FTT-IH-CHI yearly capacity "
$arr[5,4] = "ITTI"
$arr[5,5] = "TIME"
$arr[5,6] = "RSHORTTI"
$arr[5,7] = 0
$arr[5,8] = "S0"
$arr[6,0] = "IWA1"
$arr[6,1] = 1
$arr[6,2] = 0
$arr[6,3] = "FTT-IH-CHI substitution matrix"
$arr[6,4] = "ITTI"
$arr[6,5] = "ITTI"
$arr[6,6] = "NA"
$arr[6,7] = 0
$arr[6,8] = "S0"
$arr[7,0] = "BIC1"
$arr[7,1] = 1
$arr[7,2] = 0
$arr[7,3] = "FTT-IH-CHI cost matrix"
$arr[7,4] = "ITTI"
$arr[7,5] = "CTTI"
$arr[7,6] = "RSHORTTI"
$arr[7,7] = 0
$arr[7,8] = "S0"
$arr[8,0] = "IWS1"
$arr[8,1] = 1
$arr[8,2] = 0
$arr[8,3] = "FTT-IH-CHI market shares"
$arr[8,4] = "ITTI"
$arr[8,5] = "TIME"
$arr[8,6] = "RSHORTTI"
$arr[8,7] = 0
$arr[8,8] = "S0"
$arr[9,0] = "IWB1"
$arr[9,1] = 1
$arr[9,2] = 0
$arr[9,3] = "FTT-IH-CHI learning spillover matrix"
$arr[9,4] = "ITTI"
$arr[9,5] = "ITTI"
$arr[9,6] = "NA"
$arr[9,7] = 0
$arr[9,8] = "S0"
$arr[10,0] = "IFD1"
$arr[10,1] = 1
$arr[10,2] = 0
$arr[10,3] = "FTT-IH-CHI final energy demand"
$arr[10,4] = "ITTI"
$arr[10,5] = "TIME"
$arr[10,6] = "RSHORTTI"
$arr[10,7] = 0
$arr[10,8] = "S0"
$arr[11,0] = "ILC1"
$arr[11,1] = 1
$arr[11,2] = 0
$arr[11,3] = "FTT-IH-CHI The real bare LC without taxes"
$arr[11,4] = "ITTI"
$arr[11,5] = "TIME"
$arr[11,6] = "RSHORTTI"
$arr[11,7] = 0
$arr[11,8] = "S0"
$arr[12,0] = "ILG1"
$arr[12,1] = 1
$arr[12,2] = 0
$arr[12,3] = "FTT-IH-CHI LC as seen by consumer"
$arr[12,4] = "ITTI"
$arr[12,5] = "TIME"
$arr[12,6] = "RSHORTTI"
$arr[12,7] = 0
$arr[12,8] = "S0"
$arr[13,0] = "ILD1"
$arr[13,1] = 1
$arr[13,2] = 0
$arr[13,3] = "FTT-IH-CHI LC standard deviation"
$arr[13,4] = "ITTI"
$arr[13,5] = "TIME"
$arr[13,6] = "RSHORTTI"
$arr[13,7] = 0
$arr[13,8] = "S0"
$arr[14,0] = "IWE1"
$arr[14,1] = 1
$arr[14,2] = 0
$arr[14,3] = "FTT-IH-CHI Emissions"
$arr[14,4] = "ITTI"
$arr[14,5] = "TIME"
$arr[14,6] = "RSHORTTI"
$arr[14,7] = 0
$arr[14,8] = "S0"
$arr[15,0] = "IHW1"
$arr[15,1] = 1
$arr[15,2] = 0
$arr[15,3] = "FTT-IH-CHI Global average emissions per UED (kt of CO2/GWh)"
$arr[15,4] = "ITTI"
$arr[15,5] = "TIME"
$arr[15,6] = "NA"
$arr[15,7] = 0
$arr[15,8] = "S0"
$arr[16,0] = "IAM1"
$arr[16,1] = 1
$arr[16,2] = 0
$arr[16,3] = "FTT-IH-CHI gamma values"
$arr[16,4] = "ITTI"
$arr[16,5] = "TIME"
$arr[16,6] = "RSHORTTI"
$arr[16,7] = 0
$arr[16,8] = "S0"
$arr[17,0] = "IRG1"
$arr[17,1] = 1
$arr[17,2] = 0
$arr[17,3] = "FTT-IH-CHI regulations (based on capacity)"
$arr[17,4] = "ITTI"
$arr[17,5] = "TIME"
$arr[17,6] = "RSHORTTI"
$arr[17,7] = 0
$arr[17,8] = "S0"
$arr[18,0] = "ISB1"
$arr[18,1] = 1
$arr[18,2] = 0
$arr[18,3] = "FTT-IH-CHI subsidies (percentage of investment cost)"
$arr[18,4] = "ITTI"
$arr[18,5] = "TIME"
$arr[18,6] = "RSHORTTI"
$arr[18,7] = 0
$arr[18,8] = "S0"
$arr[19,0] = "IXS1"
$arr[19,1] = 1
$arr[19,2] = 0
$arr[19,3] = "FTT-IH-CHI exogenous share changes"
$arr[19,4] = "ITTI"
$arr[19,5] = "TIME"
$arr[19,6] = "RSHORTTI"
$arr[19,7] = 0
$arr[19,8] = "S0"
$arr[20,0] = "IHF1"
$arr[20,1] = 1
$arr[20,2] = 0
$arr[20,3] = "FTT-IH-CHI final fuel demand for industrial heat (ktoe)"
$arr[20,4] = "JTI"
$arr[20,5] = "TIME"
$arr[20,6] = "RSHORTTI"
$arr[20,7] = 0
$arr[20,8] = "S0"
$arr[21,0] = "IJT1"
$arr[21,1] = 1
$arr[21,2] = 0
$arr[21,3] = "FTT-IH-CHI Tech to fuel conversion matrix (fuel x technology)"
$arr[21,4] = "JTI"
$arr[21,5] = "ITTI"
$arr[21,6] = "NA"
$arr[21,7] = 0
$arr[21,8] = "S0"
$arr[22,0] = "IFT1"
$arr[22,1] = 1
$arr[22,2] = 0
$arr[22,3] = "FTT-IH-CHI Fuel tax (2010 Euros/MWh) "
$arr[22,4] = "ITTI"
$arr[22,5] = "TIME"
$arr[22,6] = "RSHORTTI"
$arr[22,7] = 0
$arr[22,8] = "S0"

$newSheet.Range("A1:I23").Value = $arr

# Restore the previously-selected cell on FTT-H (it is no longer the active
# sheet/tab once the new sheet is inserted, but keeps its own selection).
$wsH.Range("H6").Select()

# Leave the new sheet as the active tab/selection, matching the workbook
# being saved with FTT-IH-CHI on screen.
$newSheet.Activate()
$newSheet.Range("J16").Select()
